$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    # Writing a numeric-looking string straight into .Value makes Excel
    # coerce it to a real number. Instead, enter it as a text-producing
    # formula and then "Paste Special -> Values" over itself: that bakes
    # the formula result in as a literal value while preserving its Text
    # type (stored as a shared string, t="s") without touching the
    # cell's number format / style (no quote-prefix side effect).
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)
}

# Row 14: juice / 22
$ws.Cells.Item(14, 1).Value = "juice"
Set-TextValue $ws.Cells.Item(14, 2) "22"

# Row 16: pie / 22
$ws.Cells.Item(16, 1).Value = "pie"
Set-TextValue $ws.Cells.Item(16, 2) "22"

# Row 17: test / 64
$ws.Cells.Item(17, 1).Value = "test"
Set-TextValue $ws.Cells.Item(17, 2) "64"
